$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column A for the new rows as Text so date-like strings
# ("02-08-2021", etc.) are stored verbatim instead of being
# auto-converted into date serial numbers.
$ws.Range("A147:A168").NumberFormat = "@"

$ws.Range("A147").Value = "02-08-2021"
$ws.Range("B147").Value = 81041
$ws.Range("C147").Value = 14942
$ws.Range("D147").Value = 66098

$ws.Range("A148").Value = "03-08-2021"
$ws.Range("B148").Value = 81251
$ws.Range("C148").Value = 14926
$ws.Range("D148").Value = 66325

$ws.Range("A149").Value = "04-08-2021"
$ws.Range("B149").Value = 81000
$ws.Range("C149").Value = 14936
$ws.Range("D149").Value = 66065

$ws.Range("A150").Value = "05-08-2021"
$ws.Range("B150").Value = 81151
$ws.Range("C150").Value = 14971
$ws.Range("D150").Value = 66180

$ws.Range("A151").Value = "06-08-2021"
$ws.Range("B151").Value = 81017
$ws.Range("C151").Value = 15037
$ws.Range("D151").Value = 65980

$ws.Range("A152").Value = "09-08-2021"
$ws.Range("B152").Value = 81309
$ws.Range("C152").Value = 15127
$ws.Range("D152").Value = 66181

$ws.Range("A153").Value = "10-08-2021"
$ws.Range("B153").Value = 80796
$ws.Range("C153").Value = 15111
$ws.Range("D153").Value = 65685

$ws.Range("A154").Value = "11-08-2021"
$ws.Range("B154").Value = 81059
$ws.Range("C154").Value = 15090
$ws.Range("D154").Value = 65970

$ws.Range("A155").Value = "12-08-2021"
$ws.Range("B155").Value = 80948
$ws.Range("C155").Value = 15068
$ws.Range("D155").Value = 65880

$ws.Range("A156").Value = "13-08-2021"
$ws.Range("B156").Value = 80429
$ws.Range("C156").Value = 15069
$ws.Range("D156").Value = 65359

$ws.Range("A157").Value = "16-08-2021"
$ws.Range("B157").Value = 80378
$ws.Range("C157").Value = 15117
$ws.Range("D157").Value = 65261

$ws.Range("A158").Value = "17-08-2021"
$ws.Range("B158").Value = 79931
$ws.Range("C158").Value = 15084
$ws.Range("D158").Value = 64847

$ws.Range("A159").Value = "18-08-2021"
$ws.Range("B159").Value = 79462
$ws.Range("C159").Value = 15049
$ws.Range("D159").Value = 64413

$ws.Range("A160").Value = "19-08-2021"
$ws.Range("B160").Value = 79699
$ws.Range("C160").Value = 15027
$ws.Range("D160").Value = 64673

$ws.Range("A161").Value = "20-08-2021"
$ws.Range("B161").Value = 78833
$ws.Range("C161").Value = 15041
$ws.Range("D161").Value = 63792

$ws.Range("A162").Value = "23-08-2021"
$ws.Range("B162").Value = 78569
$ws.Range("C162").Value = 15063
$ws.Range("D162").Value = 63506

$ws.Range("A163").Value = "24-08-2021"
$ws.Range("B163").Value = 78923
$ws.Range("C163").Value = 15025
$ws.Range("D163").Value = 63898

$ws.Range("A164").Value = "25-08-2021"
$ws.Range("B164").Value = 79662
$ws.Range("C164").Value = 14979
$ws.Range("D164").Value = 64683

$ws.Range("A165").Value = "26-08-2021"
$ws.Range("B165").Value = 79766
$ws.Range("C165").Value = 14947
$ws.Range("D165").Value = 64819

$ws.Range("A166").Value = "27-08-2021"
$ws.Range("B166").Value = 79609
$ws.Range("C166").Value = 14934
$ws.Range("D166").Value = 64675

$ws.Range("A167").Value = "30-08-2021"
$ws.Range("B167").Value = 80797
$ws.Range("C167").Value = 14963
$ws.Range("D167").Value = 65834

$ws.Range("A168").Value = "31-08-2021"
$ws.Range("B168").Value = 80589
$ws.Range("C168").Value = 14967
$ws.Range("D168").Value = 65622

# Restore the default (Normal) cell style so the new cells match
# the unstyled data cells used throughout the rest of the sheet.
$ws.Range("A147:D168").Style = "Normal"
